# Add the missing "Emails" column to the Organizations import template.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G1").Value = "Emails"
